# Updates the cryptos worksheet with refreshed price/volume data
# (and a couple of row content swaps) as captured in the commit
# "Updated cryptos list on Sat May 13 04:12:04 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.912.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4282"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3683"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07228"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8610"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.27"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.030.89"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +8.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.621"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.380"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06896"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008852"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.950.24"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.178"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.248.81"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +9.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.66"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.883"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.207"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.887"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +15.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.85"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08940"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7423"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.161"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.430"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.797"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05215"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01923"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5083"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.732"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.95%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1641"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.429"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.249"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4580"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.652"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06282"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.803"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.78%  "
